$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (shifts existing rows 3..128 down to 4..129)
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 with its data. Columns A,B,C,E,F,G,H,I,N,Q,R hold the same
# constant values used throughout the sheet (copied from row 2); D,J,K,L,M,O,P are
# the new record's own values.
$ws.Range("A3").Value = 10
$ws.Range("B3").Value = "Vega Modelo de Temuco"
$ws.Range("C3").Value = "La Araucanía"
$ws.Range("D3").Value = 44812
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E3").Value = 9
$ws.Range("F3").Value = 100112035
$ws.Range("G3").Value = "Bruselas (repollito)"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 180
$ws.Range("K3").Value = 24000
$ws.Range("L3").Value = 25000
$ws.Range("M3").Value = 24444
$ws.Range("N3").Value = "$/malla 10 kilos"
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("P3").Value = 2444
$ws.Range("Q3").Value = 10
$ws.Range("R3").Value = "Hortaliza"
